$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 119.90909
$ws.Cells.Item(6, 9).Value = 119.90909
$ws.Cells.Item(6, 11).Value = 359.72727
$ws.Cells.Item(6, 13).Value = -247.72727
$ws.Cells.Item(62, 8).Value = 6935.4287
$ws.Cells.Item(62, 9).Value = 5244.8887
$ws.Cells.Item(62, 10).Value = 9978.4
$ws.Cells.Item(62, 11).Value = 5244.8887
$ws.Cells.Item(62, 12).Value = 9978.4
$ws.Cells.Item(62, 13).Value = -4620.8887
$ws.Cells.Item(62, 14).Value = -11226.4
$ws.Cells.Item(65, 8).Value = 6935.4287
$ws.Cells.Item(65, 9).Value = 5244.8887
$ws.Cells.Item(65, 10).Value = 9978.4
$ws.Cells.Item(65, 11).Value = 26224.4435
$ws.Cells.Item(65, 12).Value = 49892
$ws.Cells.Item(65, 13).Value = -23104.4435
$ws.Cells.Item(65, 14).Value = -56132
$ws.Cells.Item(98, 8).Value = 420.77777
$ws.Cells.Item(98, 9).Value = 404.79166
$ws.Cells.Item(98, 11).Value = 404.79166
$ws.Cells.Item(98, 13).Value = 1093.20834
$ws.Cells.Item(122, 8).Value = 420.77777
$ws.Cells.Item(122, 9).Value = 404.79166
$ws.Cells.Item(122, 11).Value = 1214.37498
$ws.Cells.Item(122, 13).Value = 1235.62502
$ws.Cells.Item(127, 8).Value = 962.3333
$ws.Cells.Item(127, 9).Value = 848.5
$ws.Cells.Item(127, 10).Value = 1019.25
$ws.Cells.Item(127, 11).Value = 2545.5
$ws.Cells.Item(127, 12).Value = 3057.75
$ws.Cells.Item(127, 13).Value = 2414.5
$ws.Cells.Item(127, 14).Value = -12977.75
$ws.Cells.Item(132, 8).Value = 29414654
$ws.Cells.Item(132, 9).Value = 32261120
$ws.Cells.Item(132, 11).Value = 96783360
$ws.Cells.Item(132, 13).Value = -96780830
$ws.Cells.Item(137, 8).Value = 86656.50999999999
$ws.Cells.Item(137, 9).Value = 122666.09
$ws.Cells.Item(137, 10).Value = 1776.7858
$ws.Cells.Item(137, 11).Value = 367998.27
$ws.Cells.Item(137, 12).Value = 5330.357400000001
$ws.Cells.Item(137, 13).Value = -365448.27
$ws.Cells.Item(137, 14).Value = -10430.3574
$ws.Cells.Item(138, 8).Value = 3314.8367
$ws.Cells.Item(138, 9).Value = 2015.5238
$ws.Cells.Item(138, 10).Value = 4289.3213
$ws.Cells.Item(138, 11).Value = 6046.5714
$ws.Cells.Item(138, 12).Value = 12867.9639
$ws.Cells.Item(138, 13).Value = -906.5713999999998
$ws.Cells.Item(138, 14).Value = -23147.9639

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1092.3334
$ws.Cells.Item(2, 9).Value = 1017.8421
$ws.Cells.Item(2, 10).Value = 1800
$ws.Cells.Item(2, 11).Value = 1017.8421
$ws.Cells.Item(2, 12).Value = 1800
$ws.Cells.Item(2, 13).Value = -904.8421
$ws.Cells.Item(2, 14).Value = -2026
$ws.Cells.Item(32, 8).Value = 9590.947
$ws.Cells.Item(32, 9).Value = 7041.1406
$ws.Cells.Item(32, 11).Value = 7041.1406
$ws.Cells.Item(32, 13).Value = -6754.1406
$ws.Cells.Item(61, 8).Value = 9806173
$ws.Cells.Item(61, 9).Value = 11906446
$ws.Cells.Item(61, 11).Value = 11906446
$ws.Cells.Item(61, 13).Value = -11906234
$ws.Cells.Item(74, 8).Value = 40001610
$ws.Cells.Item(74, 9).Value = 52632316
$ws.Cells.Item(74, 10).Value = 4366.6665
$ws.Cells.Item(74, 11).Value = 52632316
$ws.Cells.Item(74, 12).Value = 4366.6665
$ws.Cells.Item(74, 13).Value = -52631442
$ws.Cells.Item(74, 14).Value = -6114.6665
$ws.Cells.Item(77, 8).Value = 40001610
$ws.Cells.Item(77, 9).Value = 52632316
$ws.Cells.Item(77, 10).Value = 4366.6665
$ws.Cells.Item(77, 11).Value = 263161580
$ws.Cells.Item(77, 12).Value = 21833.3325
$ws.Cells.Item(77, 13).Value = -263157212
$ws.Cells.Item(77, 14).Value = -30569.3325
$ws.Cells.Item(97, 8).Value = 382.625
$ws.Cells.Item(97, 9).Value = 382.625
$ws.Cells.Item(97, 11).Value = 382.625
$ws.Cells.Item(97, 13).Value = 113.375
$ws.Cells.Item(116, 8).Value = 1092.3334
$ws.Cells.Item(116, 9).Value = 1017.8421
$ws.Cells.Item(116, 10).Value = 1800
$ws.Cells.Item(116, 11).Value = 1017.8421
$ws.Cells.Item(116, 12).Value = 1800
$ws.Cells.Item(116, 13).Value = 1276.1579
$ws.Cells.Item(116, 14).Value = -6388
$ws.Cells.Item(122, 8).Value = 2925.8125
$ws.Cells.Item(122, 9).Value = 2879.5386
$ws.Cells.Item(122, 11).Value = 8638.6158
$ws.Cells.Item(122, 13).Value = -6188.6158
$ws.Cells.Item(132, 8).Value = 9270572
$ws.Cells.Item(132, 9).Value = 11112935
$ws.Cells.Item(132, 10).Value = 58756.777
$ws.Cells.Item(132, 11).Value = 33338805
$ws.Cells.Item(132, 12).Value = 176270.331
$ws.Cells.Item(132, 13).Value = -33336275
$ws.Cells.Item(132, 14).Value = -181330.331
$ws.Cells.Item(136, 8).Value = 9806173
$ws.Cells.Item(136, 9).Value = 11906446
$ws.Cells.Item(136, 11).Value = 35719338
$ws.Cells.Item(136, 13).Value = -35716788

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1092.3334
$ws.Cells.Item(3, 9).Value = 1017.8421
$ws.Cells.Item(3, 10).Value = 1800
$ws.Cells.Item(3, 11).Value = 1017.8421
$ws.Cells.Item(3, 12).Value = 1800
$ws.Cells.Item(3, 13).Value = -903.8421
$ws.Cells.Item(3, 14).Value = -2028
$ws.Cells.Item(20, 8).Value = 2879.8
$ws.Cells.Item(20, 9).Value = 2849.75
$ws.Cells.Item(20, 11).Value = 2849.75
$ws.Cells.Item(20, 13).Value = -2602.75
$ws.Cells.Item(40, 8).Value = 34224
$ws.Cells.Item(40, 10).Value = 34224
$ws.Cells.Item(40, 12).Value = 34224
$ws.Cells.Item(40, 14).Value = -34754
$ws.Cells.Item(134, 8).Value = 4390.6943
$ws.Cells.Item(134, 9).Value = 4526.212
$ws.Cells.Item(134, 10).Value = 2900
$ws.Cells.Item(134, 11).Value = 13578.636
$ws.Cells.Item(134, 12).Value = 8700
$ws.Cells.Item(134, 13).Value = -11043.636
$ws.Cells.Item(134, 14).Value = -13770

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(97, 8).Value = 32000
$ws.Cells.Item(97, 10).Value = 32000
$ws.Cells.Item(97, 12).Value = 32000
$ws.Cells.Item(97, 14).Value = -33982
$ws.Cells.Item(99, 8).Value = 3843.3333
$ws.Cells.Item(99, 9).Value = 2882.7778
$ws.Cells.Item(99, 11).Value = 2882.7778
$ws.Cells.Item(99, 13).Value = -1384.7778
$ws.Cells.Item(126, 8).Value = 3843.3333
$ws.Cells.Item(126, 9).Value = 2882.7778
$ws.Cells.Item(126, 11).Value = 8648.3334
$ws.Cells.Item(126, 13).Value = -6178.3334
$ws.Cells.Item(132, 8).Value = 45457464
$ws.Cells.Item(132, 9).Value = 52633524
$ws.Cells.Item(132, 11).Value = 157900572
$ws.Cells.Item(132, 13).Value = -157898042
$ws.Cells.Item(134, 8).Value = 58824556
$ws.Cells.Item(134, 9).Value = 58824556
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 11).Value = 176473668
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 13).Value = -176471133
$ws.Cells.Item(134, 14).ClearContents()
$ws.Cells.Item(137, 8).Value = 21326.666
$ws.Cells.Item(137, 10).Value = 27490
$ws.Cells.Item(137, 12).Value = 27490
$ws.Cells.Item(137, 14).Value = -37690
$ws.Cells.Item(141, 8).Value = 19367.875
$ws.Cells.Item(141, 9).Value = 5000
$ws.Cells.Item(141, 10).Value = 21420.428
$ws.Cells.Item(141, 11).Value = 5000
$ws.Cells.Item(141, 12).Value = 21420.428
$ws.Cells.Item(141, 13).Value = 180
$ws.Cells.Item(141, 14).Value = -31780.428

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1278.4894
$ws.Cells.Item(5, 10).Value = 1768.9048
$ws.Cells.Item(5, 12).Value = 5306.7144
$ws.Cells.Item(5, 14).Value = -5530.7144
$ws.Cells.Item(131, 8).Value = 675.02
$ws.Cells.Item(131, 9).Value = 326.30768
$ws.Cells.Item(131, 10).Value = 727.12646
$ws.Cells.Item(131, 11).Value = 978.92304
$ws.Cells.Item(131, 12).Value = 2181.37938
$ws.Cells.Item(131, 13).Value = 4061.07696
$ws.Cells.Item(131, 14).Value = -12261.37938
$ws.Cells.Item(135, 8).Value = 1278.4894
$ws.Cells.Item(135, 10).Value = 1768.9048
$ws.Cells.Item(135, 12).Value = 15920.1432
$ws.Cells.Item(135, 14).Value = -20990.1432

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 2980502.8
$ws.Cells.Item(70, 9).Value = 4362.5835
$ws.Cells.Item(70, 10).Value = 6948690
$ws.Cells.Item(70, 11).Value = 4362.5835
$ws.Cells.Item(70, 12).Value = 6948690
$ws.Cells.Item(70, 13).Value = -4092.5835
$ws.Cells.Item(70, 14).Value = -6949230
$ws.Cells.Item(73, 8).Value = 2980502.8
$ws.Cells.Item(73, 9).Value = 4362.5835
$ws.Cells.Item(73, 10).Value = 6948690
$ws.Cells.Item(73, 11).Value = 4362.5835
$ws.Cells.Item(73, 12).Value = 6948690
$ws.Cells.Item(73, 13).Value = -3426.5835
$ws.Cells.Item(73, 14).Value = -6950562
$ws.Cells.Item(122, 8).Value = 3939.2083
$ws.Cells.Item(122, 9).Value = 3556.8948
$ws.Cells.Item(122, 10).Value = 5392
$ws.Cells.Item(122, 11).Value = 10670.6844
$ws.Cells.Item(122, 12).Value = 16176
$ws.Cells.Item(122, 13).Value = -8220.6844
$ws.Cells.Item(122, 14).Value = -21076
$ws.Cells.Item(132, 8).Value = 3754048
$ws.Cells.Item(132, 9).Value = 5525819
$ws.Cells.Item(132, 10).Value = 49436.184
$ws.Cells.Item(132, 11).Value = 16577457
$ws.Cells.Item(132, 12).Value = 148308.552
$ws.Cells.Item(132, 13).Value = -16574927
$ws.Cells.Item(132, 14).Value = -153368.552

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 1522.05
$ws.Cells.Item(22, 9).Value = 1798.6428
$ws.Cells.Item(22, 10).Value = 876.6667
$ws.Cells.Item(22, 11).Value = 1798.6428
$ws.Cells.Item(22, 12).Value = 876.6667
$ws.Cells.Item(22, 13).Value = -1503.6428
$ws.Cells.Item(22, 14).Value = -1466.6667
$ws.Cells.Item(27, 8).Value = 1522.05
$ws.Cells.Item(27, 9).Value = 1798.6428
$ws.Cells.Item(27, 10).Value = 876.6667
$ws.Cells.Item(27, 11).Value = 1798.6428
$ws.Cells.Item(27, 12).Value = 876.6667
$ws.Cells.Item(27, 13).Value = -1691.6428
$ws.Cells.Item(27, 14).Value = -1090.6667
$ws.Cells.Item(46, 8).Value = 1149.75
$ws.Cells.Item(93, 8).Value = 3352.8572
$ws.Cells.Item(93, 9).Value = 3411.6667
$ws.Cells.Item(93, 11).Value = 3411.6667
$ws.Cells.Item(93, 13).Value = -2163.6667
$ws.Cells.Item(122, 8).Value = 821692.8
$ws.Cells.Item(122, 9).Value = 1229876.8
$ws.Cells.Item(122, 10).Value = 5325
$ws.Cells.Item(122, 11).Value = 3689630.4
$ws.Cells.Item(122, 12).Value = 15975
$ws.Cells.Item(122, 13).Value = -3687180.4
$ws.Cells.Item(122, 14).Value = -20875
$ws.Cells.Item(132, 8).Value = 296938.84
$ws.Cells.Item(132, 9).Value = 448792.28
$ws.Cells.Item(132, 10).Value = 4078.5715
$ws.Cells.Item(132, 11).Value = 1346376.84
$ws.Cells.Item(132, 12).Value = 12235.7145
$ws.Cells.Item(132, 13).Value = -1343846.84
$ws.Cells.Item(132, 14).Value = -17295.7145
$ws.Cells.Item(136, 8).Value = 2248.2415
$ws.Cells.Item(136, 9).Value = 2248.2415
$ws.Cells.Item(136, 11).Value = 6744.7245
$ws.Cells.Item(136, 13).Value = -4194.7245

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(20, 8).Value = 4764.6665
$ws.Cells.Item(20, 10).Value = 5117.6
$ws.Cells.Item(20, 12).Value = 5117.6
$ws.Cells.Item(20, 14).Value = -5597.6
$ws.Cells.Item(81, 8).Value = 145
$ws.Cells.Item(81, 9).Value = 145
$ws.Cells.Item(81, 11).Value = 290
$ws.Cells.Item(81, 13).Value = 771
$ws.Cells.Item(84, 8).Value = 145
$ws.Cells.Item(84, 9).Value = 145
$ws.Cells.Item(84, 11).Value = 1450
$ws.Cells.Item(84, 13).Value = 3854
$ws.Cells.Item(136, 8).Value = 35598730
$ws.Cells.Item(136, 9).Value = 43012444
$ws.Cells.Item(136, 10).Value = 12900
$ws.Cells.Item(136, 11).Value = 129037332
$ws.Cells.Item(136, 12).Value = 38700
$ws.Cells.Item(136, 13).Value = -129034782
$ws.Cells.Item(136, 14).Value = -43800
